$wb = $excel.ActiveWorkbook

# Remove the now-unused "testDataForMethod1" sheet, keeping only
# "getConceptModelDataByCondition" as the sole worksheet.
$excel.DisplayAlerts = $false
$wsDelete = $wb.Worksheets.Item("testDataForMethod1")
$wsDelete.Delete()
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("getConceptModelDataByCondition")
$ws.Activate()

# Rename the response-related header cells to the new naming scheme.
$ws.Range("F1").Value = "rspStatus"
$ws.Range("G1").Value = "rspCode"
$ws.Range("H1").Value = "rspMessage"

# Refresh the sheet's current selection.
$ws.Range("I5").Select()
